$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Description" column header + value -------------------------------
$ws.Range("H1").Value = "Description"
$ws.Range("H2").Value = "hello brand"

# --- Quick-registration row data --------------------------------------------
$ws.Range("A2").Value = "Redline tuning"
$ws.Range("B2").Value = "Ha-Banai Street 29, Holon, 58857, Israel"
$ws.Range("C2").Value = "Event/Catering"
$ws.Range("D2").Value = "American"
$ws.Range("F2").Value = "HarriTest"
$ws.Range("G2").Value = "C:\Users\Harri\Desktop\Food.jpg"

# --- Business email becomes a live mailto hyperlink -------------------------
$ws.Range("E2").Value = "test@harri.com"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:test@harri.com") | Out-Null

# --- Selection moves to H3 ---------------------------------------------------
$ws.Range("H3").Select() | Out-Null
